$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 856; this shifts the existing rows 856-934
# down to 857-935 and carries forward formatting (incl. the date number
# format on column D) from the row that used to be at 856.
$ws.Rows.Item(856).Insert()

# Populate the newly inserted row 856 with the new data observation.
$ws.Range("A856").Value2 = 8
$ws.Range("B856").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C856").Value2 = "Coquimbo"
$ws.Range("D856").Value2 = 45132
$ws.Range("E856").Value2 = 4
$ws.Range("F856").Value2 = 100112043
$ws.Range("G856").Value2 = "Pepino ensalada"
$ws.Range("H856").Value2 = "Sin especificar"
$ws.Range("I856").Value2 = "Primera"
$ws.Range("J856").Value2 = 600
$ws.Range("K856").Value2 = 9000
$ws.Range("L856").Value2 = 10000
$ws.Range("M856").Value2 = 9500
$ws.Range("N856").Value2 = "`$/caja 60 unidades"
$ws.Range("O856").Value2 = "Región de Arica y Parinacota"
$ws.Range("P856").Value2 = 158
$ws.Range("Q856").Value2 = 60
$ws.Range("R856").Value2 = "Hortaliza"
